# Financials update - fill in the latest period (column D) values that were
# previously "NA", and correct a handful of other figures on the FLUX sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FLUX")

# Income Statement section
$ws.Range("D17").Value = 10300    # Total Operating Expenses
$ws.Range("D18").Value = -6200    # Operating Income or Loss
$ws.Range("D20").Value = 0        # Total Other Income/Expenses Net
$ws.Range("D21").Value = -6200    # Earnings Before Interest And Taxes
$ws.Range("G21").Value = -2300    # Earnings Before Interest And Taxes (2016)
$ws.Range("H21").Value = -4200    # Earnings Before Interest And Taxes (2015)
$ws.Range("J21").Value = -2400    # Earnings Before Interest And Taxes (2012)
$ws.Range("D23").Value = -7000    # Income Before Tax
$ws.Range("D26").Value = -7000    # Income After Tax
$ws.Range("D27").Value = -7000    # Net Income From Continuing Ops
$ws.Range("D32").Value = 0        # Other Items
$ws.Range("D33").Value = -7000    # Net Income
$ws.Range("D35").Value = -7000    # Net Income Applicable To Common Shares

# Balance Sheet section
$ws.Range("D58").Value = 10900    # Short/Current Long Term Debt
$ws.Range("D61").Value = 0        # Long Term Debt

# Cash Flow Statement section
$ws.Range("D81").Value = -7000    # Net Income
